$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.267433524131775
$ws.Range("B1").Value = 2.177060842514038
$ws.Range("C1").Value = 4.586565017700195
$ws.Range("D1").Value = 3.069464683532715
$ws.Range("E1").Value = 1.373181462287903
